# Updates the "cryptos" price table (Sheet1) with refreshed Price (col D)
# and Volume(1h) (col E) figures, plus a row swap: rows 46/47 (Aave /
# Stellar) exchange places along with their Coin/Link/Price/Volume data.
#
# Price values are plain text in this workbook (not numbers), e.g.
# "59.515.02" or "1.00" with a significant trailing zero. Assigning such a
# numeric-looking string straight to Range.Value would make Excel's COM
# layer "smart match" it into a real number (and round/trim it), exactly
# like typing it into a live worksheet would. To keep those cells as text
# we prefix the literal with an apostrophe (the normal Excel "force text"
# convention) and then immediately clear the resulting quote-prefix style
# back to Normal so the cell's format stays untouched, matching the
# original (unstyled) data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "59.035.05"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "2.499.23"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'537.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").Value = "'138.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.76%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").Value = "'0.560"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "2.499.51"
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("D10").Value = "'0.100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "'5.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").Value = "'0.348"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.07%  "
$ws.Range("D14").Value = "2.952.43"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "'22.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.72%  "
$ws.Range("D16").Value = "58.962.86"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").Value = "'0.0000140"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").Value = "2.501.41"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D19").Value = "'10.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("D20").Value = "'4.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").Value = "'323.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "'5.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").Value = "'62.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").Value = "'0.414"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.97%  "
$ws.Range("D26").Value = "'0.167"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "'7.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.54%  "
$ws.Range("D29").Value = "0.0₃0771"
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("D30").Value = "'6.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.32%  "
$ws.Range("D31").Value = "'1.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").Value = "'165.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "'1.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.45%  "
$ws.Range("D35").Value = "'1.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.84%  "
$ws.Range("D36").Value = "'18.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("D37").Value = "'4.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.96%  "
$ws.Range("D38").Value = "'1.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.30%  "
$ws.Range("D39").Value = "'3.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("D40").Value = "'0.799"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.76%  "
$ws.Range("D41").Value = "'5.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.33%  "
$ws.Range("D42").Value = "'278.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.54%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").Value = "'10.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "'0.593"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0935"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'124.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'0.0508"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.78%  "
$ws.Range("D49").Value = "'0.0221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.06%  "
$ws.Range("D50").Value = "'17.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.25%  "
$ws.Range("D51").Value = "1.763.61"
$ws.Range("E51").Value = "  -2.82%  "
